$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old" suffix -> "_FV2304", "_new" suffix -> "_FV2310"
#    Column layout (A1:U1): 10 "*_old" headers, "diff", 10 "*_new" headers.
$leftHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
$rightHeaders = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $leftHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $leftHeaders[$i]
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $rightHeaders.Count; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $rightHeaders[$i]
}

# 2. Turn the used range into an Excel Table ("Table1")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U91"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# 3. Freeze the header row (top row)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
